# Fix O&M window typo: refresh cumulative-capacity projections
# (Column B values recomputed; Baseline-Low/Moderate-Low lose their
#  trailing 2058 row now that the window no longer overruns.)
$wb = $excel.ActiveWorkbook

# --- Baseline-Low ---
$ws = $wb.Worksheets.Item("Baseline-Low")
$ws.Cells.Item(2, 2).Value = 1627.807136655092
$ws.Cells.Item(3, 2).Value = 2334.585451983126
$ws.Cells.Item(4, 2).Value = 3077.290613913541
$ws.Cells.Item(5, 2).Value = 3812.44496461266
$ws.Cells.Item(6, 2).Value = 4488.452146064634
$ws.Cells.Item(7, 2).Value = 5253.728661349983
$ws.Cells.Item(8, 2).Value = 6015.259188441078
$ws.Cells.Item(9, 2).Value = 6751.250087432547
$ws.Cells.Item(10, 2).Value = 7487.240986424014
$ws.Cells.Item(11, 2).Value = 8341.181536593529
$ws.Cells.Item(12, 2).Value = 9390.738377412039
$ws.Cells.Item(13, 2).Value = 10368.61519730597
$ws.Cells.Item(14, 2).Value = 11432.22740586152
$ws.Cells.Item(15, 2).Value = 12443.95682558234
$ws.Cells.Item(16, 2).Value = 13407.65323706305
$ws.Cells.Item(17, 2).Value = 14461.78947628179
$ws.Cells.Item(18, 2).Value = 15465.28774702235
$ws.Cells.Item(19, 2).Value = 16436.10267932832
$ws.Cells.Item(20, 2).Value = 17488.43301133983
$ws.Cells.Item(21, 2).Value = 18480.05251112411
$ws.Cells.Item(22, 2).Value = 19441.14378569984
$ws.Cells.Item(23, 2).Value = 20491.14363214689
$ws.Cells.Item(24, 2).Value = 21482.74501084779
$ws.Cells.Item(25, 2).Value = 22443.77372745027
$ws.Cells.Item(26, 2).Value = 23500.84207744599
$ws.Cells.Item(27, 2).Value = 24592.92580724673
$ws.Rows.Item(28).Delete()

# --- Baseline-Mid (SC) ---
$ws = $wb.Worksheets.Item("Baseline-Mid (SC)")
$ws.Cells.Item(2, 2).Value = 1014.686248331108
$ws.Cells.Item(3, 2).Value = 1723.928200563715
$ws.Cells.Item(4, 2).Value = 2652.480900977594
$ws.Cells.Item(5, 2).Value = 3737.659267944707
$ws.Cells.Item(6, 2).Value = 5405.037494671253
$ws.Cells.Item(7, 2).Value = 7286.013900313775
$ws.Cells.Item(8, 2).Value = 9144.594764214053
$ws.Cells.Item(9, 2).Value = 10954.11508666184
$ws.Cells.Item(10, 2).Value = 12768.9411359356
$ws.Cells.Item(11, 2).Value = 14623.6582912613
$ws.Cells.Item(12, 2).Value = 16559.53882147466
$ws.Cells.Item(13, 2).Value = 18470.18061703362
$ws.Cells.Item(14, 2).Value = 20290.75680477099
$ws.Cells.Item(15, 2).Value = 22268.49720168707
$ws.Cells.Item(16, 2).Value = 23376.74780385629
$ws.Cells.Item(17, 2).Value = 24496.32143330742

# --- Baseline-Mid (CC) ---
$ws = $wb.Worksheets.Item("Baseline-Mid (CC)")
$ws.Cells.Item(2, 2).Value = 1014.686248331108
$ws.Cells.Item(3, 2).Value = 1723.928200563715
$ws.Cells.Item(4, 2).Value = 2291.103519976236
$ws.Cells.Item(5, 2).Value = 3015.891875835438
$ws.Cells.Item(6, 2).Value = 3965.591786236091
$ws.Cells.Item(7, 2).Value = 5140.600315955766
$ws.Cells.Item(8, 2).Value = 6294.860774066269
$ws.Cells.Item(9, 2).Value = 7757.582049913518
$ws.Cells.Item(10, 2).Value = 9142.772446536506
$ws.Cells.Item(11, 2).Value = 10542.66790702044
$ws.Cells.Item(12, 2).Value = 11967.6337443147
$ws.Cells.Item(13, 2).Value = 13388.17813356356
$ws.Cells.Item(14, 2).Value = 14782.15159882103
$ws.Cells.Item(15, 2).Value = 16176.12506407851
$ws.Cells.Item(16, 2).Value = 17583.27467643587
$ws.Cells.Item(17, 2).Value = 19116.5632392959
$ws.Cells.Item(18, 2).Value = 20041.31234199035
$ws.Cells.Item(19, 2).Value = 20544.64490921627
$ws.Cells.Item(20, 2).Value = 21049.66753396935
$ws.Cells.Item(21, 2).Value = 21556.17230413414
$ws.Cells.Item(22, 2).Value = 22064.10408042578
$ws.Cells.Item(23, 2).Value = 22582.14074512123
$ws.Cells.Item(24, 2).Value = 23102.96673658975
$ws.Cells.Item(25, 2).Value = 23627.98921186694
$ws.Cells.Item(26, 2).Value = 24149.63952874978
$ws.Cells.Item(27, 2).Value = 24663.09126077018

# --- Moderate-Low ---
$ws = $wb.Worksheets.Item("Moderate-Low")
$ws.Cells.Item(2, 2).Value = 1919.734452711094
$ws.Cells.Item(3, 2).Value = 3062.61643202483
$ws.Cells.Item(4, 2).Value = 4245.667383870446
$ws.Cells.Item(5, 2).Value = 5430.236879244304
$ws.Cells.Item(6, 2).Value = 6552.617506566201
$ws.Cells.Item(7, 2).Value = 7748.14480555061
$ws.Cells.Item(8, 2).Value = 8938.750567751422
$ws.Cells.Item(9, 2).Value = 10737.50609513517
$ws.Cells.Item(10, 2).Value = 12560.10695459791
$ws.Cells.Item(11, 2).Value = 14502.05260734156
$ws.Cells.Item(12, 2).Value = 16677.1357781006
$ws.Cells.Item(13, 2).Value = 18863.16687652341
$ws.Cells.Item(14, 2).Value = 20901.46783272906
$ws.Cells.Item(15, 2).Value = 22443.95682558234
$ws.Cells.Item(16, 2).Value = 23407.65323706305
$ws.Cells.Item(17, 2).Value = 24461.78947628179
$ws.Cells.Item(18, 2).Value = 25465.28774702235
$ws.Cells.Item(19, 2).Value = 26436.10267932832
$ws.Cells.Item(20, 2).Value = 27488.43301133983
$ws.Cells.Item(21, 2).Value = 28480.05251112411
$ws.Cells.Item(22, 2).Value = 29441.14378569984
$ws.Cells.Item(23, 2).Value = 30491.14363214689
$ws.Cells.Item(24, 2).Value = 31482.74501084779
$ws.Cells.Item(25, 2).Value = 32443.77372745027
$ws.Cells.Item(26, 2).Value = 33500.842077446
$ws.Cells.Item(27, 2).Value = 34592.92580724673
$ws.Rows.Item(28).Delete()

# --- Moderate-Mid (SC) ---
$ws = $wb.Worksheets.Item("Moderate-Mid (SC)")
$ws.Cells.Item(2, 2).Value = 1014.686248331108
$ws.Cells.Item(3, 2).Value = 2015.855516619717
$ws.Cells.Item(4, 2).Value = 3380.511881019298
$ws.Cells.Item(5, 2).Value = 4906.036037901612
$ws.Cells.Item(6, 2).Value = 6662.083208256897
$ws.Cells.Item(7, 2).Value = 8628.686858723344
$ws.Cells.Item(8, 2).Value = 10559.36558821764
$ws.Cells.Item(9, 2).Value = 12408.56885251923
$ws.Cells.Item(10, 2).Value = 14808.69772646339
$ws.Cells.Item(11, 2).Value = 17232.67193248653
$ws.Cells.Item(12, 2).Value = 19660.32550125241
$ws.Cells.Item(13, 2).Value = 22120.83530867213
$ws.Cells.Item(14, 2).Value = 24664.06972337513
$ws.Cells.Item(15, 2).Value = 27025.60835442773
$ws.Cells.Item(16, 2).Value = 29028.31723864588
$ws.Cells.Item(17, 2).Value = 30613.57506482587
$ws.Cells.Item(18, 2).Value = 31574.98171177761
$ws.Cells.Item(19, 2).Value = 32033.85980095197
$ws.Cells.Item(20, 2).Value = 32508.97879705755
$ws.Cells.Item(21, 2).Value = 32982.79965382951
$ws.Cells.Item(22, 2).Value = 33465.58570483124
$ws.Cells.Item(23, 2).Value = 33948.70946393117
$ws.Cells.Item(24, 2).Value = 34445.76877234803
$ws.Cells.Item(25, 2).Value = 34942.95930529542

# --- Expanded-High ---
$ws = $wb.Worksheets.Item("Expanded-High")
$ws.Cells.Item(2, 2).Value = 1014.686248331108
$ws.Cells.Item(3, 2).Value = 1723.928200563715
$ws.Cells.Item(4, 2).Value = 3437.346855886964
$ws.Cells.Item(5, 2).Value = 5305.246735263695
$ws.Cells.Item(6, 2).Value = 7409.072479854041
$ws.Cells.Item(7, 2).Value = 10147.91602383643
$ws.Cells.Item(8, 2).Value = 12964.92538126792
$ws.Cells.Item(9, 2).Value = 16027.06846302629
$ws.Cells.Item(10, 2).Value = 19469.91728751813
$ws.Cells.Item(11, 2).Value = 23302.33724353648
$ws.Cells.Item(12, 2).Value = 27506.59834337482
$ws.Cells.Item(13, 2).Value = 32004.84905238706
$ws.Cells.Item(14, 2).Value = 36447.59295925179
$ws.Cells.Item(15, 2).Value = 40904.52398475842
$ws.Cells.Item(16, 2).Value = 45300.97908368556
$ws.Cells.Item(17, 2).Value = 49297.02427482433
$ws.Cells.Item(18, 2).Value = 52622.46709841703
$ws.Cells.Item(19, 2).Value = 54572.09618688664
